# Regenerate s_val data (filter save games) for bruihl_justin 2023 sheet.
# Updates columns B (TB), C (d2S), D (K), E (IP) and G (sum) for data rows 2-10.
# Column F (Win) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729 }
    3  = @{ B = 0.3464964993005633;  C = 0.004309184025731883; D = 0.1529057820181812;  E = 0.4998867070740569; G = 1.003598172418533 }
    4  = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729 }
    5  = @{ B = 1.505614041169197;   C = 1.65323645889881;    D = 3.082599426703578;   E = 0.4998867070740569; G = 6.741336633845642 }
    6  = @{ B = 0.7287194209349384;  C = 1.65323645889881;    D = 0.1529057820181812;  E = 6.48142807727062;   G = 9.016289739122548 }
    7  = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729 }
    8  = @{ B = 0.006876353814593728; C = 0.05231270169004087; D = 16.98373111632243;  E = 6.48142807727062;   G = 23.52434824909768 }
    9  = @{ B = 0.7287194209349384;  C = 1.65323645889881;    D = 3.082599426703578;   E = 0.4998867070740569; G = 5.964442013611383 }
    10 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
